$wb = $excel.ActiveWorkbook

# --- About sheet: label the workbook with the state name ---------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "Colorado"

# --- BECbIC sheet: refresh the "Employee Compensation" row with the ----
# --- recalculated OECD-derived values -----------------------------------
$wsBec = $wb.Worksheets.Item("BECbIC")

$newValues = [ordered]@{
    "B2"  = 102097641.0212372
    "C2"  = 145798886.6971826
    "D2"  = 1513881662.826288
    "E2"  = 268371049.7329473
    "F2"  = 1443671642.80235
    "G2"  = 1746587092.129685
    "H2"  = 84328864.69514453
    "I2"  = 224404879.5371597
    "J2"  = 421165952.6602732
    "K2"  = 140456000
    "L2"  = 281937445.1684638
    "M2"  = 289860069.4348322
    "N2"  = 286766353.7149698
    "O2"  = 73876274.04260629
    "P2"  = 474256324.2100542
    "Q2"  = 0
    "R2"  = 32509126.38744986
    "S2"  = 932165173.309472
    "T2"  = 1253431981.793264
    "U2"  = 177451716.7693134
    "V2"  = 1124843144.354261
    "W2"  = 128582559.6067797
    "X2"  = 0
    "Y2"  = 798506358.5665708
    "Z2"  = 991004683.8900487
    "AA2" = 269598259.9164927
    "AB2" = 584837295.7550453
    "AC2" = 13472408753.72844
    "AD2" = 10478084929.19887
    "AE2" = 6699010000
    "AF2" = 8204030606.930657
    "AG2" = 8423663078.411325
    "AH2" = 4274308000
    "AI2" = 2341004000
    "AJ2" = 14659013000
    "AK2" = 3310616078.759262
    "AL2" = 30610286797.09698
    "AM2" = 38782018000
    "AN2" = 2657194000
    "AO2" = 20092488000
    "AP2" = 2559927968.24097
    "AQ2" = 0
}

foreach ($addr in $newValues.Keys) {
    $wsBec.Range($addr).Value = $newValues[$addr]
}
